# Dev new report: split the single "Tháng 2,3,4" sheet into "BLHH" (the
# existing report) and a new "SXCN" sheet that mirrors the same data.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Rename the existing sheet -------------------------------------------------
$ws1.Name = "BLHH"

# --- Add the new sheet right after BLHH and name it -----------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "SXCN"

# --- Duplicate BLHH's data + formatting into SXCN --------------------------------
$ws1.Range("A1:M5").Copy($ws2.Range("A1:M5"))

# Column widths on SXCN differ from BLHH for a handful of columns (others keep
# the sheet default width). COM ColumnWidth is expressed in "characters" and is
# offset from the raw stored width by the standard 5/6 character padding.
$offset = 0.8333333333333334
$ws2.Columns.Item(1).ColumnWidth = (0 - $offset)
$ws2.Columns.Item(3).ColumnWidth = (19.88671875 - $offset)
$ws2.Columns.Item(7).ColumnWidth = (17.77734375 - $offset)
$ws2.Columns.Item(8).ColumnWidth = (18.21875 - $offset)
$ws2.Columns.Item(9).ColumnWidth = (22.77734375 - $offset)
$ws2.Columns.Item(10).ColumnWidth = (21.21875 - $offset)
$ws2.Columns.Item(11).ColumnWidth = (22.5546875 - $offset)
$ws2.Columns.Item(12).ColumnWidth = (19.77734375 - $offset)
$ws2.Columns.Item(13).ColumnWidth = (17.88671875 - $offset)

# Row heights on SXCN were re-computed by Excel (not manually dragged) after the
# paste, because the narrower default-width columns wrap the text differently.
$ws2.Rows.Item(1).RowHeight = 100.8
$ws2.Rows.Item(2).RowHeight = 57.6
$ws2.Rows.Item(3).RowHeight = 57.6
$ws2.Rows.Item(4).RowHeight = 115.2
$ws2.Rows.Item(5).RowHeight = 86.4

# --- View / selection state -----------------------------------------------------
# BLHH is no longer the active tab; its saved selection now spans the whole
# used range.
$ws1.Range("A1:M5").Select()

# SXCN becomes the active sheet/tab, with a simple single-cell selection.
$ws2.Activate()
$ws2.Range("D3").Select()
